# AtlasOfLivingAustralia/biocollect#1759 - fixes bug with empty table throwing an error
#
# Adds three new rows (8, 9, 10) of sample data to the bottom of the
# "form1" worksheet's bulk-import example table, and updates the
# selected cell to reflect where the user ended up after the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 8: serial = 3, a = "test a 3", b.c = "test bc 6"
$ws.Cells.Item(8, 1).Value = 3
$ws.Cells.Item(8, 2).Value = "test a 3"
$ws.Cells.Item(8, 3).Value = "test bc 6"

# New row 9: serial = 3 only (rest of the row intentionally blank)
$ws.Cells.Item(9, 1).Value = 3

# New row 10: serial = 3 only (rest of the row intentionally blank)
$ws.Cells.Item(10, 1).Value = 3

# Move the active selection to L8, and scroll the sheet back so column A
# is the leftmost visible column again.
[void]$ws.Range("L8").Select()
